# Auto-generated Excel COM-interop script to apply cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.138.49"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "1.911.24"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7416"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "245.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3099"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.47"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06989"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08082"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7711"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").Value = "1.927.93"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.354"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "30.132.13"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.980"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007870"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.93%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.160.15"
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.134"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.426"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1293"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.063"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.71%  "
$ws.Range("E30").Value = "  +2.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.351"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.346"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.102"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05167"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7517"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.736"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01951"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.801"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.353"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4521"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.996"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8420"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.807"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.946"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("D49").Value = "2.060.27"
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1187"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.73%  "
